# Week 6 - Day 1 - Assignment 1
#
# Target shape of the workbook after this script runs:
#   sheets (in tab order): names, companies, leads, contact
#     - "companies" is the original "contact" sheet, renamed, with a new
#       "Company" column (TCS / Wipro / IQVIA) added next to the Phone list.
#     - "leads" keeps its data except row 3 col A: "Azim" -> "pavi".
#     - "contact" is a brand-new sheet (added at the end) holding the same
#       Phone / 97 / 98 / 87 list the old "contact" sheet used to have in
#       column A, and ends up as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- Rename the existing "contact" sheet to "companies" ------------------
$companies = $wb.Worksheets.Item("contact")
$companies.Name = "companies"

# --- Add a fresh sheet at the end; this becomes the new "contact" sheet --
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$contact = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$contact.Name = "contact"

# Populate the new "contact" sheet with the Phone list (text, not numbers -
# leading apostrophe keeps "97"/"98"/"87" stored as strings like before).
$contact.Range("A1").Value = "Phone"
$contact.Range("A1").Font.Bold = $true
$contact.Range("A2").Value = "'97"
$contact.Range("A3").Value = "'98"
$contact.Range("A4").Value = "'87"

# --- Add the "Company" column to the "companies" sheet -------------------
$companies.Range("B1").Value = "Company"
$companies.Range("B1").Font.Bold = $true
$companies.Range("B2").Value = "TCS"
$companies.Range("B3").Value = "Wipro"
$companies.Range("B4").Value = "IQVIA"
$companies.Columns.Item(2).ColumnWidth = 9.7
[void]$companies.Range("E15").Select()

# --- Fix the "leads" sheet: "Azim" -> "pavi" ------------------------------
$leads = $wb.Worksheets.Item("leads")
$leads.Range("A3").Value = "pavi"
[void]$leads.Range("A3").Select()

# --- Leave the new "contact" sheet selected/active, like the source file -
[void]$contact.Activate()
[void]$contact.Range("D14").Select()
